$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.868.77"
$ws.Range("E2").Value = "  +0.91%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.286.52"
$ws.Range("E3").Value = "  +2.24%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.10"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.642"
$ws.Range("E6").Value = "  +2.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.55"
$ws.Range("E7").Value = "  +5.96%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.645"
$ws.Range("E9").Value = "  +2.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.05"
$ws.Range("E10").Value = "  -2.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0978"
$ws.Range("E11").Value = "  +2.89%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "59.00"
$ws.Range("E12").Value = "  -1.00%  "
$ws.Range("E13").Value = "  +3.09%  "
$ws.Range("E14").Value = "  +1.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.627.82"
$ws.Range("E15").Value = "  +2.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.28"
$ws.Range("E16").Value = "  +3.81%  "
$ws.Range("E17").Value = "  -1.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.285.67"
$ws.Range("E18").Value = "  +2.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.799.60"
$ws.Range("E19").Value = "  +0.98%  "
$ws.Range("E20").Value = "  +3.46%  "
$ws.Range("E21").Value = "  +1.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.48"
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.69"
$ws.Range("E23").Value = "  +1.92%  "
$ws.Range("E24").Value = "  +6.59%  "
$ws.Range("E25").Value = "  -0.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.52"
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("E28").Value = "  -0.62%  "
$ws.Range("E29").Value = "  -1.07%  "
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "166.73"
$ws.Range("E31").Value = "  -0.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.02"
$ws.Range("E32").Value = "  +1.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.47"
$ws.Range("E33").Value = "  +6.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.127"
$ws.Range("E34").Value = "  +3.48%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0823"
$ws.Range("E35").Value = "  +4.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.09"
$ws.Range("E36").Value = "  +10.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.127"
$ws.Range("E37").Value = "  +2.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.61"
$ws.Range("E38").Value = "  +11.12%  "
$ws.Range("E39").Value = "  +1.77%  "
$ws.Range("E40").Value = "  -3.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "14.44"
$ws.Range("E41").Value = "  +13.60%  "
$ws.Range("E42").Value = "  +3.10%  "
$ws.Range("E43").Value = "  +3.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.218"
$ws.Range("E44").Value = "  +8.90%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.78"
$ws.Range("E45").Value = "  -1.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.07"
$ws.Range("E46").Value = "  +4.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.89"
$ws.Range("E47").Value = "  -3.67%  "
$ws.Range("E48").Value = "  +2.22%  "
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "99.97"
$ws.Range("E51").Value = "  +6.20%  "
